# Update "Horarios Línea 141" workbook with the 01:56:31 scrape results.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01:56:31"
$ws1.Range("A3").Value = "Total filas: 5"

# Existing row 7 gets refreshed with the new scrape time / minutes.
$ws1.Range("A7").Value = "01:56:31"
$ws1.Range("D7").Value = 2

# Two brand-new rows appended for this scrape.
$ws1.Range("A9").Value = "01:56:31"
$ws1.Range("B9").Value = "02:59"
$ws1.Range("C9").Value = "215_ALUAR"
$ws1.Range("D9").Value = 63
$ws1.Range("E9").Value = "LP1912"

$ws1.Range("A10").Value = "01:56:31"
$ws1.Range("B10").Value = "03:48"
$ws1.Range("C10").Value = "14_ABASTO"
$ws1.Range("D10").Value = 112
$ws1.Range("E10").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 01:56:31"
$ws2.Range("A3").Value = "Total filas: 3"

# One brand-new row appended for this scrape.
$ws2.Range("A8").Value = "01:56:31"
$ws2.Range("B8").Value = "02:59"
$ws2.Range("C8").Value = "215_ALUAR"
$ws2.Range("D8").Value = 63
$ws2.Range("E8").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 01:56:31"
